# Commit: "add language id for property name"
#
# Row 9 of the "Property1" sheet holds, per-column, a human-readable
# description of the property in column 1's Id/Name. These descriptions
# (previously hard-coded Chinese text) are replaced by a language id key
# ("LPID_<ENUM_NAME>") that mirrors the property's enum name in row 1,
# so the client can look the text up from a localization table instead.
#
# Row 1 (B1:AD1) already holds the bare enum names (SUCKBLOOD,
# REFLECTDAMAGE, CRITICAL, MAXHP, ...) - row 9 gets "LPID_" + that name.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

$lpidNames = @(
    "LPID_SUCKBLOOD",
    "LPID_REFLECTDAMAGE",
    "LPID_CRITICAL",
    "LPID_MAXHP",
    "LPID_MAXMP",
    "LPID_MAXSP",
    "LPID_HPREGEN",
    "LPID_SPREGEN",
    "LPID_MPREGEN",
    "LPID_ATK_VALUE",
    "LPID_DEF_VALUE",
    "LPID_MOVE_SPEED",
    "LPID_ATK_SPEED",
    "LPID_ATK_FIRE",
    "LPID_ATK_LIGHT",
    "LPID_ATK_WIND",
    "LPID_ATK_ICE",
    "LPID_ATK_POISON",
    "LPID_DEF_FIRE",
    "LPID_DEF_LIGHT",
    "LPID_DEF_WIND",
    "LPID_DEF_ICE",
    "LPID_DEF_POISON",
    "LPID_DIZZY_GATE",
    "LPID_MOVE_GATE",
    "LPID_SKILL_GATE",
    "LPID_PHYSICAL_GATE",
    "LPID_MAGIC_GATE",
    "LPID_BUFF_GATE"
)

# Row 9, columns B (2) .. AD (30) - one LPID per property column.
$descRow = 9
$firstCol = 2

for ($i = 0; $i -lt $lpidNames.Length; $i++) {
    $col = $firstCol + $i
    $ws.Cells.Item($descRow, $col).Value = $lpidNames[$i]
}

# Keep selection on the last edited cell (AD9), matching the authoring
# session that made this change.
$ws.Range("AD9").Select()
